$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Time Lists" actual hours value (I3) from 56 to 60
$ws.Range("I3").Value = 60

# Update the active selection to K8 (matches the saved selection state in the diff)
$ws.Range("K8").Select()
